$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = '@'
    $c.Value = $text
    $c.Style = 'Normal'
}

Set-TextValue 'D2' '63.663.33'
Set-TextValue 'E2' '  +0.21%  '
Set-TextValue 'D3' '2.621.76'
Set-TextValue 'E3' '  -0.79%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '595.19'
Set-TextValue 'E5' '  -1.40%  '
Set-TextValue 'D6' '150.45'
Set-TextValue 'E6' '  +2.86%  '
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'D8' '0.588'
Set-TextValue 'E8' '  +0.00%  '
Set-TextValue 'E9' '  +0.32%  '
Set-TextValue 'D10' '5.71'
Set-TextValue 'E10' '  +2.35%  '
Set-TextValue 'E11' '  +2.97%  '
Set-TextValue 'E12' '  -1.17%  '
Set-TextValue 'D13' '27.72'
Set-TextValue 'E13' '  +0.56%  '
Set-TextValue 'D14' '3.090.17'
Set-TextValue 'E14' '  -0.87%  '
Set-TextValue 'D15' '63.449.45'
Set-TextValue 'E15' '  +0.18%  '
Set-TextValue 'D16' '0.0000152'
Set-TextValue 'E16' '  +3.52%  '
Set-TextValue 'D17' '2.611.98'
Set-TextValue 'E17' '  -2.15%  '
Set-TextValue 'D18' '12.34'
Set-TextValue 'E18' '  +7.47%  '
Set-TextValue 'D19' '4.66'
Set-TextValue 'E19' '  +1.71%  '
Set-TextValue 'D20' '347.28'
Set-TextValue 'E20' '  +0.93%  '
Set-TextValue 'D21' '6.87'
Set-TextValue 'E21' '  -0.99%  '
Set-TextValue 'D22' '0.998'
Set-TextValue 'E22' '  -0.18%  '
Set-TextValue 'E23' '  +2.73%  '
Set-TextValue 'D24' '66.33'
Set-TextValue 'E24' '  -0.36%  '
Set-TextValue 'D25' '1.73'
Set-TextValue 'E25' '  +11.58%  '
Set-TextValue 'D26' '9.25'
Set-TextValue 'E26' '  +1.39%  '
Set-TextValue 'D27' '1.68'
Set-TextValue 'E27' '  -0.80%  '
Set-TextValue 'D28' '567.65'
Set-TextValue 'E28' '  -2.05%  '
Set-TextValue 'D29' '8.25'
Set-TextValue 'E29' '  +3.94%  '
Set-TextValue 'D30' '0.161'
Set-TextValue 'E30' '  -1.09%  '
Set-TextValue 'E31' '  +0.05%  '
Set-TextValue 'D33' '0.0₃0847'
Set-TextValue 'E33' '  +2.51%  '
Set-TextValue 'D34' '1.76'
Set-TextValue 'E34' '  -0.30%  '
Set-TextValue 'D35' '5.25'
Set-TextValue 'E35' '  +0.53%  '
Set-TextValue 'D36' '168.51'
Set-TextValue 'E36' '  +0.98%  '
Set-TextValue 'D38' '0.999'
Set-TextValue 'E38' '  -0.05%  '
Set-TextValue 'E39' '  -0.46%  '
Set-TextValue 'D40' '19.39'
Set-TextValue 'E40' '  +1.41%  '
Set-TextValue 'E41' '  +0.06%  '
Set-TextValue 'D42' '167.00'
Set-TextValue 'D43' '39.90'
Set-TextValue 'E43' '  -0.20%  '
Set-TextValue 'D44' '3.92'
Set-TextValue 'E44' '  +3.86%  '
Set-TextValue 'D45' '0.0597'
Set-TextValue 'E45' '  +4.85%  '
Set-TextValue 'D46' '21.54'
Set-TextValue 'E46' '  -2.64%  '
Set-TextValue 'D47' '0.629'
Set-TextValue 'E47' '  -0.24%  '
Set-TextValue 'D48' '0.0250'
Set-TextValue 'E48' '  +1.47%  '
Set-TextValue 'D49' '1.99'
Set-TextValue 'E49' '  +4.72%  '
Set-TextValue 'D50' '0.0963'
Set-TextValue 'E50' '  +0.13%  '
Set-TextValue 'D51' '19.35'
Set-TextValue 'E51' '  +2.96%  '
